$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet, row 3 (the b2545090 row): update Correspond Handoff/Handback Datetime
$wsZhCn.Range("D3").Value = "2016-01-18 02:16:43"
$wsZhCn.Range("G3").Value = "2016-01-18 02:17:32"

# de-de sheet, row 3 (the b2545090 row): update Correspond Handoff/Handback Datetime
$wsDeDe.Range("D3").Value = "2016-01-18 02:16:55"
$wsDeDe.Range("G3").Value = "2016-01-18 02:17:52"
